$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update LOC_METHOD, CYCLO_METHOD, WMC_CLASS, NOM_CLASS values for row 2
$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 2.0
$ws.Range("G2").Value = 5.0
$ws.Range("H2").Value = 1.0
